$wb = $excel.ActiveWorkbook

# The "Variables" sheet has a row for the FIGURES pseudo-variable whose
# variable-code (column C) is "persons". It previously also had an explicit
# variable-label (en/da/kl, columns E/F/G) of "antal"/"antal"/"amerlassusaat".
# Clear those labels so the variable-code is used as the fallback label.
$ws = $wb.Worksheets.Item("Variables")
$ws.Range("E5:G5").ClearContents()

# Make "Variables" the active/selected sheet and select E5:F5 (matching the
# recorded selection after the edit was made in Excel).
$ws.Activate()
$ws.Range("E5:F5").Select()
